$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header text updates (January -> May) ---
$ws.Range("C2").Value = "Turbo English - May"
$ws.Range("J2").Value = "May"
$ws.Range("J3").Value = "5 Classes"

# --- Date row (row 6) updates ---
$ws.Range("A6").Value = "05/05"
$ws.Range("C6").Value = "05/06"
$ws.Range("E6").Value = "05/07"
$ws.Range("G6").Value = "05/08"
$ws.Range("I6").Value = "05/09"

# --- Extend the title merge to include column H ---
$ws.Range("C2:H3").Merge()
$ws.Range("D2:H2").Style = "Normal"
$ws.Range("C3:H3").Style = "Normal"

# --- Row 7 content: first re-style G7/I7 (currently the red "No School" style)
#     to match the plain A7/C7/E7 style, since that style is retired in the
#     target workbook. Copy formatting only (values are overwritten next). ---
$ws.Range("A7").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("I7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 7 new text content ---
$ws.Range("A7").Value = "Welome:`nWelcome`n4 - 7"
$ws.Range("B7").Value = "Workbook:`np. 4 - 5`n Study flashcards "
$ws.Range("C7").Value = "Unit  1:`nIn the Classroom`n8 - 11"
$ws.Range("D7").Value = "Workbook:`np. 6 - 7`n Study flashcards "
$ws.Range("E7").Value = "Unit  1:`nIn the Classroom`n12 - 15"
$ws.Range("F7").Value = "Workbook:`np. 8 - 9`n Study flashcards "
$ws.Range("G7").Value = "Test"
$ws.Range("H7").Value = "Flashcards"
$ws.Range("I7").Value = "Review"
$ws.Range("J7").Value = "Flashcards"
